# Commit: Mon, Jul 27, 2020  3:05:31 AM
#
# 1) Slide 5's table switches to a different table style (new tableStyleId GUID).
# 2) The deck's theme colours change from the "Integral" (Red Violet) palette
#    to the standard "Office" palette (the palette that previously only lived
#    in the notes-master's theme part). The font scheme / format scheme are
#    already identical between the two theme parts, so only the colour scheme
#    needs to move.

$p = $ppt.ActivePresentation

# --- 1. Table style on slide 5 --------------------------------------------
$slide = $p.Slides.Item(5)
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shape = $slide.Shapes.Item($i)
    if ($shape.HasTable) {
        $shape.Table.ApplyStyle("{6A998459-FBD4-43F6-9FEF-1863ACB112F9}", $true)
    }
}

# --- 2. Swap the presentation's theme colour scheme for the Office palette -
$cs = $p.SlideMaster.ColorScheme
$cs.Colors(1).RGB  = 0          # dk1      000000
$cs.Colors(2).RGB  = 16777215   # lt1      FFFFFF
$cs.Colors(3).RGB  = 6968388    # dk2      44546A
$cs.Colors(4).RGB  = 15132391   # lt2      E7E6E6
$cs.Colors(5).RGB  = 13998939   # accent1  5B9BD5
$cs.Colors(6).RGB  = 3243501    # accent2  ED7D31
$cs.Colors(7).RGB  = 10855845   # accent3  A5A5A5
$cs.Colors(8).RGB  = 49407      # accent4  FFC000
$cs.Colors(9).RGB  = 12874308   # accent5  4472C4
$cs.Colors(10).RGB = 4697456    # accent6  70AD47
$cs.Colors(11).RGB = 12673797   # hlink    0563C1
$cs.Colors(12).RGB = 7491477    # folHlink 954F72
